$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Re-order the "Periodo Mora" column (E16:E19) from descending (2307..2304) to
# ascending (2304..2307); the "Valor Mora" (F) stays tied to its row, so swap
# the values that moved between row 16 and row 19 accordingly.
$ws.Range("E16").Value = "2304"
$ws.Range("E17").Value = "2305"
$ws.Range("E18").Value = "2306"
$ws.Range("E19").Value = "2307"

$ws.Range("F16").Value = 42000
$ws.Range("F19").Value = 22000
